$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrdenSalida")

# Column AB ("NroOrdenSalida") holds text-typed order numbers. Rows 2-13
# previously all read "5265"; the new codelco output order number is "8762".
# A leading apostrophe keeps the cell text-typed (matching the existing
# text-stored values) instead of Excel auto-converting the digit string to
# a number.
for ($row = 2; $row -le 13; $row++) {
    $ws.Range("AB$row").Value = "'8762"
}
